$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.731.79'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.602.96'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.64'
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  +0.95%  '

$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").Value = '1.828.62'

$ws.Range("D13").Value = '1.606.55'
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.02'
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '209.25'
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("E20").Value = '  +1.53%  '

$ws.Range("E21").Value = '  +0.27%  '

$ws.Range("E22").Value = '  -5.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.05'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.61'
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.36'
$ws.Range("E28").Value = '  +0.44%  '

$ws.Range("E29").Value = '  -0.89%  '

$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.27'
$ws.Range("E31").Value = '  +1.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.96'
$ws.Range("E32").Value = '  +0.79%  '

$ws.Range("D33").Value = '1.286.81'
$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("E34").Value = '  +1.09%  '

$ws.Range("E35").Value = '  +19.15%  '

$ws.Range("E36").Value = '  +0.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.588'
$ws.Range("E37").Value = '  -5.05%  '

$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.826'
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.780'
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.78'
$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("D44").Value = '1.739.96'
$ws.Range("E44").Value = '  +0.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.37'
$ws.Range("E45").Value = '  -0.76%  '

$ws.Range("E46").Value = '  +0.21%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("E47").Value = '  +0.93%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0513'
$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.82'
$ws.Range("E49").Value = '  +13.69%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.54'
$ws.Range("E50").Value = '  +2.28%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.09%  '
